$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- update status text "Ready for handoff" -> "Handback transform failed" everywhere it occurs
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: set Error Detail (P3) and widen column P
$wsZhCn.Range("P3").Value = "Handback file name: lexpmcdu.peu is different with handoff file name: ed3c665b-0787-419c-8328-392049ddabfc.1e99952a8d6c29ef8d1a78d870b757115f267d43.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: set Error Detail (P3) and widen column P
$wsDeDe.Range("P3").Value = "Handback file name: lexpmcdu.peu is different with handoff file name: ed3c665b-0787-419c-8328-392049ddabfc.1e99952a8d6c29ef8d1a78d870b757115f267d43.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
